$d = $word.ActiveDocument

# Paragraph 3 is the "Add brick animation ... be played." bullet (split across
# several runs because of the embedded grammar-check proofErr markers).
# Replace its text with the mushroom-animation sentence, but stop one
# character short of Range.End so the trailing paragraph mark (and the
# bookmark that lives after the last run) is left untouched.
$p = $d.Paragraphs.Item(3)
$r = $p.Range
$r.End = $r.End - 1
$r.Text = "Create mushroom animation once Mario hits question mark."

# The following paragraph (originally paragraph 4) duplicated that same
# sentence - remove it entirely now that its text lives in paragraph 3.
$dup = $d.Paragraphs.Item(4)
$dup.Range.Delete()
